# PerformData.xlsx edit: header row rework ("Changed ms to sec")
#
# The sheet's header row (row 1) tracked "test student" / "browser" columns
# alongside the Date column; it is reworked so the header tracks
# Student / Date / Build / Domain instead, with a new Student column
# inserted at the very front (column A). The per-step timing columns
# (step 1..step 8, 9a, 9b, 9c) are untouched. The stray "IE" sample value
# that had been left in row 2 (under the old "browser" column) is cleared
# since that column no longer holds browser values.
#
# Order of the text writes below matters only for the shared-string table
# layout on save (new unique strings are appended in first-write order) —
# it doesn't change any cell's final value/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: was "test student" -> now "Date"
$ws.Range("C1").Value = "Date"

# Column D: was "Date" -> now "Build"
$ws.Range("D1").Value = "Build"

# Column E: was "browser" -> now "Domain"
$ws.Range("E1").Value = "Domain"

# New column A header
$ws.Range("A1").Value = "Student"

# Clear the leftover "IE" sample value from row 2 (old browser column)
$ws.Range("E2").Value = ""

# Move the active selection to B8 (matches the saved sheet view state)
[void]$ws.Range("B8").Select()
